$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks a price-history table: columns A..BP hold one timestamped
# price snapshot per scrape, followed by "nom" (product name) and
# "url_produit" (product URL). This commit adds one more timestamped
# snapshot column (2026-01-30 19:24:20) right before "nom"/"url_produit",
# which pushes those two columns one slot to the right (BQ->BR, BR->BS).

# Insert a new blank column at BQ; everything from BQ onward (nom, url_produit,
# and all their data) shifts right by one column automatically.
$ws.Columns("BQ:BQ").Insert()

# Give the new header cell the same formatting as the preceding timestamp
# header cell (bold, centered, bordered), then set its value.
$ws.Range("BP1").Copy() | Out-Null
$ws.Range("BQ1").PasteSpecial(-4122) | Out-Null
$ws.Range("BQ1").Value = "2026-01-30 19:24:20"

# For every product row that already has a price in the previous snapshot
# column (BP), carry that same price forward into the new BQ snapshot column
# (rows 2-80). Rows 81-206 have no price recorded in BP (blank cells), so
# their new BQ cell is left blank too - which is already the case after the
# column insert.
$lastDataRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastDataRow; $r++) {
    $prevPrice = $ws.Cells.Item($r, 68).Value2
    if ($prevPrice -ne $null -and $prevPrice -ne "") {
        $ws.Cells.Item($r, 69).Value = $prevPrice
    }
}
